$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.490.19"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").Value = "'2.457.22"
$ws.Range("E3").Value = "  -1.28%  "
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "  +0.73%  "
$ws.Range("D5").Value = "'310.87"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "'90.42"
$ws.Range("E6").Value = "  -5.31%  "
$ws.Range("E7").Value = "  -4.10%  "
$ws.Range("E8").Value = "  +0.57%  "
$ws.Range("D9").Value = "'0.485"
$ws.Range("E9").Value = "  -5.14%  "
$ws.Range("D10").Value = "'31.90"
$ws.Range("E10").Value = "  -6.65%  "
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("D13").Value = "'2.835.40"
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("D14").Value = "'6.72"
$ws.Range("E14").Value = "  -4.50%  "
$ws.Range("D15").Value = "'15.14"
$ws.Range("E15").Value = "  +2.15%  "
$ws.Range("D16").Value = "'2.437.18"
$ws.Range("E16").Value = "  -2.46%  "
$ws.Range("D17").Value = "'0.758"
$ws.Range("E17").Value = "  -4.48%  "
$ws.Range("D18").Value = "'41.272.56"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("D19").Value = "'6.18"
$ws.Range("E19").Value = "  -3.99%  "
$ws.Range("E20").Value = "  -1.85%  "
$ws.Range("D21").Value = "'69.47"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "'10.84"
$ws.Range("E22").Value = "  -7.79%  "
$ws.Range("D23").Value = "'231.84"
$ws.Range("E23").Value = "  -2.47%  "
$ws.Range("E24").Value = "  -4.67%  "
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  -4.80%  "
$ws.Range("D27").Value = "'23.72"
$ws.Range("E27").Value = "  -4.66%  "
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("D29").Value = "'9.51"
$ws.Range("E29").Value = "  -2.76%  "
$ws.Range("D30").Value = "'35.26"
$ws.Range("E30").Value = "  -3.94%  "
$ws.Range("D31").Value = "'151.35"
$ws.Range("E31").Value = "  -2.31%  "
$ws.Range("D32").Value = "'5.29"
$ws.Range("E32").Value = "  -6.54%  "
$ws.Range("E33").Value = "  -3.41%  "
$ws.Range("D34").Value = "'0.0747"
$ws.Range("E34").Value = "  -1.76%  "
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").Value = "'17.63"
$ws.Range("E35").Value = "  +1.92%  "
$ws.Range("B36").Value = "ApeXProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D36").Value = "'2.48"
$ws.Range("E36").Value = "  -3.23%  "
$ws.Range("D37").Value = "'2.90"
$ws.Range("E37").Value = "  -5.01%  "
$ws.Range("E38").Value = "  -5.69%  "
$ws.Range("E39").Value = "  -3.13%  "
$ws.Range("D40").Value = "'0.0992"
$ws.Range("E40").Value = "  -7.38%  "
$ws.Range("D41").Value = "'3.99"
$ws.Range("E41").Value = "  -0.97%  "
$ws.Range("E42").Value = "  +1.01%  "
$ws.Range("D43").Value = "'18.95"
$ws.Range("E43").Value = "  -10.98%  "
$ws.Range("D44").Value = "'1.932.27"
$ws.Range("E44").Value = "  -3.40%  "
$ws.Range("D45").Value = "'0.0276"
$ws.Range("E45").Value = "  -4.01%  "
$ws.Range("E46").Value = "  -6.90%  "
$ws.Range("D47").Value = "'8.57"
$ws.Range("E47").Value = "  -1.90%  "
$ws.Range("D48").Value = "'2.682.59"
$ws.Range("E48").Value = "  -1.60%  "
$ws.Range("D49").Value = "'94.10"
$ws.Range("E49").Value = "  -4.35%  "
$ws.Range("E50").Value = "  -5.48%  "
$ws.Range("D51").Value = "'65.18"
$ws.Range("E51").Value = "  -6.71%  "
